$d = $word.ActiveDocument

$d.Content.Find.Execute("74×78=", $true, $false, $false, $false, $false, $true, 1, $false, "98×99=", 2) | Out-Null
$d.Content.Find.Execute("77×22=", $true, $false, $false, $false, $false, $true, 1, $false, "19×95=", 2) | Out-Null
$d.Content.Find.Execute("42×46=", $true, $false, $false, $false, $false, $true, 1, $false, "65×49=", 2) | Out-Null
$d.Content.Find.Execute("27×96=", $true, $false, $false, $false, $false, $true, 1, $false, "17×50=", 2) | Out-Null
$d.Content.Find.Execute("20×50=", $true, $false, $false, $false, $false, $true, 1, $false, "78×32=", 2) | Out-Null
$d.Content.Find.Execute("24×48=", $true, $false, $false, $false, $false, $true, 1, $false, "28×58=", 2) | Out-Null
$d.Content.Find.Execute("72×82=", $true, $false, $false, $false, $false, $true, 1, $false, "22×18=", 2) | Out-Null
$d.Content.Find.Execute("44×53=", $true, $false, $false, $false, $false, $true, 1, $false, "26×26=", 2) | Out-Null
$d.Content.Find.Execute("41×35=", $true, $false, $false, $false, $false, $true, 1, $false, "55×23=", 2) | Out-Null
$d.Content.Find.Execute("99×63=", $true, $false, $false, $false, $false, $true, 1, $false, "22×45=", 2) | Out-Null
$d.Content.Find.Execute("31×71=", $true, $false, $false, $false, $false, $true, 1, $false, "84×52=", 2) | Out-Null
$d.Content.Find.Execute("26×16=", $true, $false, $false, $false, $false, $true, 1, $false, "13×93=", 2) | Out-Null
$d.Content.Find.Execute("24×19=", $true, $false, $false, $false, $false, $true, 1, $false, "71×68=", 2) | Out-Null
$d.Content.Find.Execute("31×38=", $true, $false, $false, $false, $false, $true, 1, $false, "76×79=", 2) | Out-Null
$d.Content.Find.Execute("15×74=", $true, $false, $false, $false, $false, $true, 1, $false, "64×23=", 2) | Out-Null
$d.Content.Find.Execute("51×88=", $true, $false, $false, $false, $false, $true, 1, $false, "47×68=", 2) | Out-Null
$d.Content.Find.Execute("59×42=", $true, $false, $false, $false, $false, $true, 1, $false, "38×21=", 2) | Out-Null
$d.Content.Find.Execute("80×92=", $true, $false, $false, $false, $false, $true, 1, $false, "25×87=", 2) | Out-Null
$d.Content.Find.Execute("14×53=", $true, $false, $false, $false, $false, $true, 1, $false, "21×82=", 2) | Out-Null
$d.Content.Find.Execute("47×23=", $true, $false, $false, $false, $false, $true, 1, $false, "63×17=", 2) | Out-Null
$d.Content.Find.Execute("81×81=", $true, $false, $false, $false, $false, $true, 1, $false, "68×28=", 2) | Out-Null
$d.Content.Find.Execute("41×41=", $true, $false, $false, $false, $false, $true, 1, $false, "39×62=", 2) | Out-Null
$d.Content.Find.Execute("40×83=", $true, $false, $false, $false, $false, $true, 1, $false, "13×13=", 2) | Out-Null
$d.Content.Find.Execute("95×73=", $true, $false, $false, $false, $false, $true, 1, $false, "36×23=", 2) | Out-Null
$d.Content.Find.Execute("95×89=", $true, $false, $false, $false, $false, $true, 1, $false, "88×82=", 2) | Out-Null
